# Updating odds data for Jogos_da_Semana_FlashScore_2025-01-29.xlsx
# Applies the numeric value changes described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("Q2").Value = 2.42
$ws.Range("R2").Value = 1.57
$ws.Range("U2").Value = 6
$ws.Range("G3").Value = 2.2
$ws.Range("I3").Value = 3.9
$ws.Range("J3").Value = 3
$ws.Range("M3").Value = 1.11
$ws.Range("N3").Value = 6.5
$ws.Range("V3").Value = 1.2
$ws.Range("AG3").Value = 21
$ws.Range("AJ3").Value = 6
$ws.Range("AM3").Value = 8.5
$ws.Range("G4").Value = 3.1
$ws.Range("I4").Value = 2.5
$ws.Range("J4").Value = 4
$ws.Range("L4").Value = 3.5
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("V4").Value = 1.16
$ws.Range("AF4").Value = 34
$ws.Range("AK4").Value = 21
$ws.Range("AN4").Value = 10
$ws.Range("AP4").Value = 26
$ws.Range("AQ4").Value = 29
$ws.Range("G5").Value = 2.88
$ws.Range("H5").Value = 2.88
$ws.Range("I5").Value = 2.8
$ws.Range("L5").Value = 3.75
$ws.Range("M5").Value = 1.13
$ws.Range("N5").Value = 6
$ws.Range("AA5").Value = 2.25
$ws.Range("AB5").Value = 1.57
$ws.Range("AN5").Value = 12
$ws.Range("G6").Value = 1.42
$ws.Range("H6").Value = 4.75
$ws.Range("I6").Value = 7
$ws.Range("J6").Value = 1.95
$ws.Range("K6").Value = 2.25
$ws.Range("L6").Value = 7.5
$ws.Range("Q6").Value = 1.5
$ws.Range("R6").Value = 2.6
$ws.Range("V6").Value = 1.45
$ws.Range("AD6").Value = 6
$ws.Range("AE6").Value = 8.5
$ws.Range("AF6").Value = 9
$ws.Range("AL6").Value = 81
$ws.Range("AN6").Value = 34
$ws.Range("AO6").Value = 21
$ws.Range("AP6").Value = 81
$ws.Range("AQ6").Value = 51
$ws.Range("AR6").Value = 51
$ws.Range("K13").Value = 2.05
$ws.Range("O13").Value = 1.36
$ws.Range("P13").Value = 3
$ws.Range("S13").Value = 2.2
$ws.Range("T13").Value = 1.65
$ws.Range("W13").Value = 4
$ws.Range("X13").Value = 1.22
$ws.Range("AI13").Value = 8.5
$ws.Range("G14").Value = 4
$ws.Range("H14").Value = 3.7
$ws.Range("I14").Value = 1.85
$ws.Range("L14").Value = 2.5
$ws.Range("M14").Value = 1.05
$ws.Range("N14").Value = 11
$ws.Range("O14").Value = 1.29
$ws.Range("P14").Value = 3.5
$ws.Range("S14").Value = 1.98
$ws.Range("T14").Value = 1.88
$ws.Range("AA14").Value = 1.91
$ws.Range("AB14").Value = 1.91
$ws.Range("AD14").Value = 21
$ws.Range("AK14").Value = 17
$ws.Range("AN14").Value = 8.5
$ws.Range("AP14").Value = 15
$ws.Range("G17").Value = 2.3
$ws.Range("H17").Value = 2.8
$ws.Range("I17").Value = 3.6
$ws.Range("J17").Value = 3.2
$ws.Range("L17").Value = 4.33
$ws.Range("O17").Value = 1.53
$ws.Range("P17").Value = 2.38
$ws.Range("Q17").Value = 2.1
$ws.Range("R17").Value = 1.78
$ws.Range("S17").Value = 2.7
$ws.Range("T17").Value = 1.44
$ws.Range("W17").Value = 5.5
$ws.Range("X17").Value = 1.14
$ws.Range("Y17").Value = 1.62
$ws.Range("Z17").Value = 2.2
$ws.Range("AA17").Value = 2.2
$ws.Range("AB17").Value = 1.62
$ws.Range("AD17").Value = 9.5
$ws.Range("AF17").Value = 21
$ws.Range("AK17").Value = 19
$ws.Range("AL17").Value = 81
$ws.Range("AP17").Value = 41
$ws.Range("AR17").Value = 51
